$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44294; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=500; K=7000; L=8000; M=7500; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=125; Q=60; R='Hortaliza'}
    @{Row=3; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44371; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=300; K=8500; L=9000; M=8750; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=146; Q=60; R='Hortaliza'}
    @{Row=4; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44503; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=1100; K=6500; L=7000; M=6750; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=112; Q=60; R='Hortaliza'}
    @{Row=5; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44490; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=600; K=13000; L=15000; M=14000; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=233; Q=60; R='Hortaliza'}
    @{Row=6; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44258; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=500; K=7000; L=8000; M=7500; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=125; Q=60; R='Hortaliza'}
    @{Row=7; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44377; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=400; K=7000; L=8000; M=7500; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=125; Q=60; R='Hortaliza'}
    @{Row=8; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44314; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=1100; K=7000; L=8000; M=7500; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=125; Q=60; R='Hortaliza'}
    @{Row=9; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44266; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=600; K=6500; L=7000; M=6750; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=112; Q=60; R='Hortaliza'}
    @{Row=10; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44286; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=600; K=7000; L=8000; M=7500; N='$/caja 50 unidades'; O='Provincia de Limarí'; P=150; Q=50; R='Hortaliza'}
    @{Row=11; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44335; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=500; K=7500; L=8000; M=7750; N='$/caja 50 unidades'; O='Provincia de Limarí'; P=155; Q=50; R='Hortaliza'}
    @{Row=12; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44497; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=800; K=7500; L=8000; M=7750; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=129; Q=60; R='Hortaliza'}
    @{Row=13; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44482; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=400; K=11000; L=12000; M=11500; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=192; Q=60; R='Hortaliza'}
    @{Row=14; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44328; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=500; K=7500; L=8000; M=7750; N='$/caja 50 unidades'; O='Provincia de Limarí'; P=155; Q=50; R='Hortaliza'}
    @{Row=15; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44203; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=300; K=4500; L=5000; M=4750; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=79; Q=60; R='Hortaliza'}
    @{Row=16; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44483; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=340; K=10000; L=11000; M=10500; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=175; Q=60; R='Hortaliza'}
    @{Row=17; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44217; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=700; K=6500; L=7000; M=6750; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=112; Q=60; R='Hortaliza'}
    @{Row=18; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44244; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=500; K=5000; L=6000; M=5500; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=92; Q=60; R='Hortaliza'}
    @{Row=19; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44293; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=400; K=7000; L=8000; M=7500; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=125; Q=60; R='Hortaliza'}
    @{Row=20; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44308; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=400; K=6000; L=7000; M=6500; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=108; Q=60; R='Hortaliza'}
    @{Row=21; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44321; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=500; K=7000; L=8000; M=7500; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=125; Q=60; R='Hortaliza'}
    @{Row=22; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44265; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=500; K=6500; L=7000; M=6750; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=112; Q=60; R='Hortaliza'}
    @{Row=23; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44300; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=400; K=6000; L=7000; M=6500; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=108; Q=60; R='Hortaliza'}
    @{Row=24; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44216; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=1100; K=5500; L=6000; M=5750; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=96; Q=60; R='Hortaliza'}
    @{Row=25; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44336; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=600; K=8500; L=9000; M=8750; N='$/caja 50 unidades'; O='Provincia de Limarí'; P=175; Q=50; R='Hortaliza'}
    @{Row=26; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44301; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=300; K=6000; L=7000; M=6500; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=108; Q=60; R='Hortaliza'}
    @{Row=27; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44279; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=500; K=7000; L=8000; M=7500; N='$/caja 50 unidades'; O='Provincia de Limarí'; P=150; Q=50; R='Hortaliza'}
    @{Row=28; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44504; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=700; K=6500; L=7000; M=6750; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=112; Q=60; R='Hortaliza'}
    @{Row=29; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44384; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=300; K=7000; L=8000; M=7500; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=125; Q=60; R='Hortaliza'}
    @{Row=30; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44315; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=500; K=7000; L=8000; M=7500; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=125; Q=60; R='Hortaliza'}
    @{Row=31; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44510; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=900; K=5000; L=6000; M=5500; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=92; Q=60; R='Hortaliza'}
    @{Row=32; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44517; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=500; K=5000; L=6000; M=5500; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=92; Q=60; R='Hortaliza'}
    @{Row=33; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44238; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=400; K=7000; L=8000; M=7500; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=125; Q=60; R='Hortaliza'}
    @{Row=34; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44251; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=700; K=6500; L=7000; M=6750; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=112; Q=60; R='Hortaliza'}
    @{Row=35; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44181; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=900; K=4500; L=5000; M=4750; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=79; Q=60; R='Hortaliza'}
    @{Row=36; A=2; B='Comercializadora del Agro de Limarí'; C='Coquimbo'; D=44307; E=4; F=100112032; G='Zapallo italiano'; H='Sin especificar'; I='Primera'; J=700; K=6000; L=7000; M=6500; N='$/caja 60 unidades'; O='Provincia de Limarí'; P=108; Q=60; R='Hortaliza'}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
    $ws.Cells.Item($r, 8).Value = $item.H
    $ws.Cells.Item($r, 9).Value = $item.I
    $ws.Cells.Item($r, 10).Value = $item.J
    $ws.Cells.Item($r, 11).Value = $item.K
    $ws.Cells.Item($r, 12).Value = $item.L
    $ws.Cells.Item($r, 13).Value = $item.M
    $ws.Cells.Item($r, 14).Value = $item.N
    $ws.Cells.Item($r, 15).Value = $item.O
    $ws.Cells.Item($r, 16).Value = $item.P
    $ws.Cells.Item($r, 17).Value = $item.Q
    $ws.Cells.Item($r, 18).Value = $item.R
}
